# Q3 Update - 2025
# Applies the quarterly refresh to the UN-IRN "fromCSV" sheet:
#   1. The generated short-url (column B) was rotated for the whole sheet:
#      "c4Glt5" -> "nMwKn4".
#   2. Row 183 (Afghanistan -> Iran, 2024) got revised figures:
#        refugees (N183)          : 3752317 -> 3477082
#        returned_refugees (P183) : 348      -> 237452
#   3. The trailing rows for Kuwait/Pakistan/Uzbekistan (rows 185-187) were
#      dropped from the export, shrinking the sheet from A1:V187 to A1:V184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rotate the short-url value across the whole "short-url" column ---
$ws.Columns("B").Replace("c4Glt5", "nMwKn4")

# --- 2. Revise the two figures on row 183 ---
# The source data stores every value as text, so force a text number format
# before writing the digits (otherwise Excel auto-detects them as numbers),
# then restore the original cell style (copied from an untouched neighbour
# on the same row) so only the value - not the look - of the cell changes.
$ws.Range("N183").NumberFormat = "@"
$ws.Range("N183").Value = "3477082"

$ws.Range("P183").NumberFormat = "@"
$ws.Range("P183").Value = "237452"

$ws.Range("O183").Copy()
$ws.Range("N183:P183").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Drop the three trailing rows that were removed from the export ---
$ws.Rows("185:187").Delete()
